# Auto-generated edit script applying the Goblin_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H, I, J, K, L, M, N)
# across the ALC, ARM, BSM, CUL, GSM, LTW, WVR sheets (CRP sheet has no changes).

$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H17").Value = 1041.762
$ws_ALC.Range("J17").Value = 1062.561
$ws_ALC.Range("L17").Value = 3187.683
$ws_ALC.Range("N17").Value = -3523.683
$ws_ALC.Range("H62").Value = 59375.25
$ws_ALC.Range("I62").Value = 137917.33
$ws_ALC.Range("J62").Value = 12250
$ws_ALC.Range("K62").Value = 137917.33
$ws_ALC.Range("L62").Value = 12250
$ws_ALC.Range("M62").Value = -137293.33
$ws_ALC.Range("N62").Value = -13498
$ws_ALC.Range("H65").Value = 59375.25
$ws_ALC.Range("I65").Value = 137917.33
$ws_ALC.Range("J65").Value = 12250
$ws_ALC.Range("K65").Value = 689586.6499999999
$ws_ALC.Range("L65").Value = 61250
$ws_ALC.Range("M65").Value = -686466.6499999999
$ws_ALC.Range("N65").Value = -67490
$ws_ALC.Range("H76").Value = 4926.857
$ws_ALC.Range("I76").Value = 4872.5
$ws_ALC.Range("K76").Value = 4872.5
$ws_ALC.Range("M76").Value = -4557.5
$ws_ALC.Range("H79").Value = 4926.857
$ws_ALC.Range("I79").Value = 4872.5
$ws_ALC.Range("K79").Value = 4872.5
$ws_ALC.Range("M79").Value = -3780.5
$ws_ALC.Range("H86").Value = 4447.5
$ws_ALC.Range("I86").Value = 4429
$ws_ALC.Range("J86").Value = 4490.6665
$ws_ALC.Range("K86").Value = 4429
$ws_ALC.Range("L86").Value = 4490.6665
$ws_ALC.Range("M86").Value = -3306
$ws_ALC.Range("N86").Value = -6736.6665
$ws_ALC.Range("H88").Value = 3336.4119
$ws_ALC.Range("J88").Value = 3920.7856
$ws_ALC.Range("L88").Value = 3920.7856
$ws_ALC.Range("N88").Value = -4732.7856
$ws_ALC.Range("H89").Value = 4447.5
$ws_ALC.Range("I89").Value = 4429
$ws_ALC.Range("J89").Value = 4490.6665
$ws_ALC.Range("K89").Value = 22145
$ws_ALC.Range("L89").Value = 22453.3325
$ws_ALC.Range("M89").Value = -16529
$ws_ALC.Range("N89").Value = -33685.3325
$ws_ALC.Range("H91").Value = 3336.4119
$ws_ALC.Range("J91").Value = 3920.7856
$ws_ALC.Range("L91").Value = 3920.7856
$ws_ALC.Range("N91").Value = -6728.7856
$ws_ALC.Range("H100").Value = 7129.4443
$ws_ALC.Range("I100").Value = 5500
$ws_ALC.Range("J100").Value = 7333.125
$ws_ALC.Range("K100").Value = 5500
$ws_ALC.Range("L100").Value = 7333.125
$ws_ALC.Range("M100").Value = -4959
$ws_ALC.Range("N100").Value = -8415.125
$ws_ALC.Range("H108").Value = 52500
$ws_ALC.Range("J108").Value = 52500
$ws_ALC.Range("L108").Value = 52500
$ws_ALC.Range("N108").Value = -60180
$ws_ALC.Range("H116").Value = 4674.222
$ws_ALC.Range("J116").Value = 5009.7144
$ws_ALC.Range("L116").Value = 5009.7144
$ws_ALC.Range("N116").Value = -11893.7144
$ws_ALC.Range("H129").Value = 1741.1538
$ws_ALC.Range("I129").Value = 848.3333
$ws_ALC.Range("J129").Value = 3750
$ws_ALC.Range("K129").Value = 2544.9999
$ws_ALC.Range("L129").Value = 11250
$ws_ALC.Range("M129").Value = 2455.0001
$ws_ALC.Range("N129").Value = -21250
$ws_ALC.Range("H132").Value = 1719.8108
$ws_ALC.Range("I132").Value = 1396.2188
$ws_ALC.Range("J132").Value = 3790.8
$ws_ALC.Range("K132").Value = 4188.6564
$ws_ALC.Range("L132").Value = 11372.4
$ws_ALC.Range("M132").Value = -1658.6564
$ws_ALC.Range("N132").Value = -16432.4
$ws_ALC.Range("H137").Value = 19248.666
$ws_ALC.Range("I137").Value = 35329.668
$ws_ALC.Range("J137").Value = 3167.6667
$ws_ALC.Range("K137").Value = 105989.004
$ws_ALC.Range("L137").Value = 9503.000100000001
$ws_ALC.Range("M137").Value = -103439.004
$ws_ALC.Range("N137").Value = -14603.0001

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H32").Value = 5260
$ws_ARM.Range("I32").Value = 5260
$ws_ARM.Range("K32").Value = 5260
$ws_ARM.Range("M32").Value = -4973
$ws_ARM.Range("H45").Value = 2097.2727
$ws_ARM.Range("I45").Value = 1581.5714
$ws_ARM.Range("J45").Value = 2999.75
$ws_ARM.Range("K45").Value = 1581.5714
$ws_ARM.Range("L45").Value = 2999.75
$ws_ARM.Range("M45").Value = -1204.5714
$ws_ARM.Range("N45").Value = -3753.75
$ws_ARM.Range("H92").Value = 48516.668
$ws_ARM.Range("J92").Value = 48516.668
$ws_ARM.Range("L92").Value = 48516.668
$ws_ARM.Range("N92").Value = -53508.668
$ws_ARM.Range("H110").Value = 988.8461
$ws_ARM.Range("I110").Value = 988.8461
$ws_ARM.Range("K110").Value = 988.8461
$ws_ARM.Range("M110").Value = 1056.1539

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H6").Value = 50000
$ws_BSM.Range("J6").Value = 50000
$ws_BSM.Range("L6").Value = 50000
$ws_BSM.Range("N6").Value = -50226
$ws_BSM.Range("H20").Value = 1069
$ws_BSM.Range("I20").Value = 954
$ws_BSM.Range("K20").Value = 954
$ws_BSM.Range("M20").Value = -707
$ws_BSM.Range("H86").Value = 3414.3076
$ws_BSM.Range("I86").Value = 3133.2222
$ws_BSM.Range("J86").Value = 4046.75
$ws_BSM.Range("K86").Value = 3133.2222
$ws_BSM.Range("L86").Value = 4046.75
$ws_BSM.Range("M86").Value = -2010.2222
$ws_BSM.Range("N86").Value = -6292.75
$ws_BSM.Range("H89").Value = 3414.3076
$ws_BSM.Range("I89").Value = 3133.2222
$ws_BSM.Range("J89").Value = 4046.75
$ws_BSM.Range("K89").Value = 15666.111
$ws_BSM.Range("L89").Value = 20233.75
$ws_BSM.Range("M89").Value = -10050.111
$ws_BSM.Range("N89").Value = -31465.75
$ws_BSM.Range("H105").Value = 2106.6667
$ws_BSM.Range("I105").Value = 2198.3076
$ws_BSM.Range("K105").Value = 2198.3076
$ws_BSM.Range("M105").Value = -451.3076000000001

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H5").Value = 940.35297
$ws_CUL.Range("I5").Value = 387.44446
$ws_CUL.Range("J5").Value = 1562.375
$ws_CUL.Range("K5").Value = 1162.33338
$ws_CUL.Range("L5").Value = 4687.125
$ws_CUL.Range("M5").Value = -1050.33338
$ws_CUL.Range("N5").Value = -4911.125
$ws_CUL.Range("H11").Value = 32150.094
$ws_CUL.Range("I11").Value = 40788.32
$ws_CUL.Range("K11").Value = 122364.96
$ws_CUL.Range("M11").Value = -122224.96
$ws_CUL.Range("H49").Value = 1199.5
$ws_CUL.Range("J49").Value = 999
$ws_CUL.Range("L49").Value = 2997
$ws_CUL.Range("N49").Value = -3309
$ws_CUL.Range("H98").Value = 1436.5
$ws_CUL.Range("I98").Value = 997
$ws_CUL.Range("J98").Value = 1876
$ws_CUL.Range("K98").Value = 2991
$ws_CUL.Range("L98").Value = 5628
$ws_CUL.Range("M98").Value = -1493
$ws_CUL.Range("N98").Value = -8624
$ws_CUL.Range("H113").Value = 1494.5
$ws_CUL.Range("J113").Value = 1753.7778
$ws_CUL.Range("L113").Value = 5261.3334
$ws_CUL.Range("N113").Value = -9601.3334
$ws_CUL.Range("H135").Value = 940.35297
$ws_CUL.Range("I135").Value = 387.44446
$ws_CUL.Range("J135").Value = 1562.375
$ws_CUL.Range("K135").Value = 3487.00014
$ws_CUL.Range("L135").Value = 14061.375
$ws_CUL.Range("M135").Value = -952.0001400000001
$ws_CUL.Range("N135").Value = -19131.375
$ws_CUL.Range("H140").Value = 63942.875
$ws_CUL.Range("I140").Value = 72476.14
$ws_CUL.Range("K140").Value = 217428.42
$ws_CUL.Range("M140").Value = -212248.42

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H2").Value = 981.1818
$ws_GSM.Range("I2").Value = 1709.5834
$ws_GSM.Range("K2").Value = 1709.5834
$ws_GSM.Range("M2").Value = -1596.5834
$ws_GSM.Range("H33").Value = 0
$ws_GSM.Range("J33").Value = 0
$ws_GSM.Range("L33").Value = 0
$ws_GSM.Range("N33").ClearContents()
$ws_GSM.Range("H80").Value = 7380.2915
$ws_GSM.Range("I80").Value = 11275.417
$ws_GSM.Range("K80").Value = 11275.417
$ws_GSM.Range("M80").Value = -10277.417
$ws_GSM.Range("H83").Value = 7380.2915
$ws_GSM.Range("I83").Value = 11275.417
$ws_GSM.Range("K83").Value = 56377.085
$ws_GSM.Range("M83").Value = -51385.085
$ws_GSM.Range("H96").Value = 50001
$ws_GSM.Range("J96").Value = 50001
$ws_GSM.Range("L96").Value = 50001
$ws_GSM.Range("N96").Value = -55493
$ws_GSM.Range("H97").Value = 5417.6
$ws_GSM.Range("I97").Value = 462.06668
$ws_GSM.Range("J97").Value = 20284.2
$ws_GSM.Range("K97").Value = 462.06668
$ws_GSM.Range("L97").Value = 20284.2
$ws_GSM.Range("M97").Value = 33.93331999999998
$ws_GSM.Range("N97").Value = -21276.2
$ws_GSM.Range("H107").Value = 665.8333
$ws_GSM.Range("I107").Value = 559
$ws_GSM.Range("K107").Value = 559
$ws_GSM.Range("M107").Value = 1361

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H7").Value = 4411
$ws_LTW.Range("I7").Value = 4628.4287
$ws_LTW.Range("J7").Value = 3650
$ws_LTW.Range("K7").Value = 4628.4287
$ws_LTW.Range("L7").Value = 3650
$ws_LTW.Range("M7").Value = -4516.4287
$ws_LTW.Range("N7").Value = -3874
$ws_LTW.Range("H87").Value = 50000
$ws_LTW.Range("J87").Value = 50000
$ws_LTW.Range("L87").Value = 50000
$ws_LTW.Range("N87").Value = -52246
$ws_LTW.Range("H90").Value = 50000
$ws_LTW.Range("J90").Value = 50000
$ws_LTW.Range("L90").Value = 150000
$ws_LTW.Range("N90").Value = -161232
$ws_LTW.Range("H126").Value = 4411
$ws_LTW.Range("I126").Value = 4628.4287
$ws_LTW.Range("J126").Value = 3650
$ws_LTW.Range("K126").Value = 13885.2861
$ws_LTW.Range("L126").Value = 10950
$ws_LTW.Range("M126").Value = -11415.2861
$ws_LTW.Range("N126").Value = -15890
$ws_LTW.Range("H132").Value = 4714.304
$ws_LTW.Range("I132").Value = 4394.5713
$ws_LTW.Range("J132").Value = 5211.6665
$ws_LTW.Range("K132").Value = 13183.7139
$ws_LTW.Range("L132").Value = 15634.9995
$ws_LTW.Range("M132").Value = -10653.7139
$ws_LTW.Range("N132").Value = -20694.9995
$ws_LTW.Range("H133").Value = 0
$ws_LTW.Range("J133").Value = 0
$ws_LTW.Range("L133").Value = 0
$ws_LTW.Range("N133").ClearContents()

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H62").Value = 8000
$ws_WVR.Range("I62").Value = 5000
$ws_WVR.Range("K62").Value = 5000
$ws_WVR.Range("M62").Value = -4376
$ws_WVR.Range("H65").Value = 8000
$ws_WVR.Range("I65").Value = 5000
$ws_WVR.Range("K65").Value = 25000
$ws_WVR.Range("M65").Value = -21880
$ws_WVR.Range("H95").Value = 21817.25
$ws_WVR.Range("J95").Value = 21817.25
$ws_WVR.Range("L95").Value = 21817.25
$ws_WVR.Range("N95").Value = -27309.25
$ws_WVR.Range("H119").Value = 114333
$ws_WVR.Range("J119").Value = 114333
$ws_WVR.Range("L119").Value = 114333
$ws_WVR.Range("N119").Value = -124009
$ws_WVR.Range("H132").Value = 4762.9565
$ws_WVR.Range("J132").Value = 0
$ws_WVR.Range("L132").Value = 0
$ws_WVR.Range("N132").ClearContents()
